$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Primary (default) footer -> image2.png becomes image1.png ---
$footerPrimary = $sec.Footers.Item(1)
if ($footerPrimary.Exists -and $footerPrimary.Range.InlineShapes.Count -ge 1) {
    $pic = $footerPrimary.Range.InlineShapes.Item(1)
    if ($pic.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $shp = $pic.ConvertToShape()
        $shp.Name = "image1.png"
        $shp.ConvertToInlineShape() | Out-Null
    }
}

# --- First-page footer -> image2.png becomes image1.png ---
$footerFirst = $sec.Footers.Item(2)
if ($footerFirst.Exists -and $footerFirst.Range.InlineShapes.Count -ge 1) {
    $pic = $footerFirst.Range.InlineShapes.Item(1)
    if ($pic.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $shp = $pic.ConvertToShape()
        $shp.Name = "image1.png"
        $shp.ConvertToInlineShape() | Out-Null
    }
}

# --- First-page header -> image1.jpg becomes image2.jpg ---
$headerFirst = $sec.Headers.Item(2)
if ($headerFirst.Exists -and $headerFirst.Range.InlineShapes.Count -ge 1) {
    $pic = $headerFirst.Range.InlineShapes.Item(1)
    if ($pic.AlternativeText -eq "BTec_Logo-Orange") {
        $shp = $pic.ConvertToShape()
        $shp.Name = "image2.jpg"
        $shp.ConvertToInlineShape() | Out-Null
    }
}
